# Insert a new weekly record for "Vega Monumental Concepción" / Espinaca.
# A brand-new row is inserted at row 125 (pushing the former rows 125-133
# down to 126-134), and is populated with the new week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 125..133 down to 126..134, creating a blank row 125.
$ws.Rows.Item(125).Insert()

# Populate the newly inserted row 125 with this week's values.
$ws.Cells.Item(125, 1).Value = 11
$ws.Cells.Item(125, 2).Value = 'Vega Monumental Concepción'
$ws.Cells.Item(125, 3).Value = 'Bíobío'
$ws.Cells.Item(125, 4).Value = 45147
$ws.Cells.Item(125, 5).Value = 8
$ws.Cells.Item(125, 6).Value = 100112012
$ws.Cells.Item(125, 7).Value = 'Espinaca'
$ws.Cells.Item(125, 8).Value = 'Sin especificar'
$ws.Cells.Item(125, 9).Value = 'Primera'
$ws.Cells.Item(125, 10).Value = 50
$ws.Cells.Item(125, 11).Value = 6000
$ws.Cells.Item(125, 12).Value = 6500
$ws.Cells.Item(125, 13).Value = 6200
$ws.Cells.Item(125, 14).Value = '$/cuna 10 kilos'
$ws.Cells.Item(125, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(125, 16).Value = 620
$ws.Cells.Item(125, 17).Value = 10
$ws.Cells.Item(125, 18).Value = 'Hortaliza'

# Match the date cell's number format used by the rest of column D.
$ws.Cells.Item(125, 4).NumberFormat = 'YYYY-MM-DD HH:MM:SS'
